# Add new game rows (43-49) scraped from Data/bombay1.xlsx sourced games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43
$ws.Range("A43").Value = "2024-04-27 18:43:58"
$ws.Range("B43").Value = 77
$ws.Range("C43").Value = 26
$ws.Range("D43").Value = 6
$ws.Range("E43").Value = 10
$ws.Range("F43").Value = 10
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0.001
$ws.Range("J43").Value = 0.05
$ws.Range("K43").Value = 0.003
$ws.Range("L43").Value = 100
$ws.Range("M43").Value = 500
$ws.Range("N43").Value = 10
$ws.Range("O43").Value = 9
$ws.Range("P43").Value = 2
$ws.Range("Q43").Value = 500
$ws.Range("R43").Value = 3
$ws.Range("S43").Value = 1
$ws.Range("T43").Value = 20
$ws.Range("U43").Value = 0.3376623376623377
$ws.Range("V43").Value = "Data/bombay1.xlsx"
$ws.Range("W43").Value = -206500
$ws.Range("X43").Value = "No es Simulación"

# Row 44
$ws.Range("A44").Value = "2024-04-27 18:59:38"
$ws.Range("B44").Value = 61
$ws.Range("C44").Value = 51
$ws.Range("D44").Value = 14
$ws.Range("E44").Value = 37
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0.001
$ws.Range("J44").Value = 0.05
$ws.Range("K44").Value = 0.003
$ws.Range("L44").Value = 100
$ws.Range("M44").Value = 500
$ws.Range("N44").Value = 10
$ws.Range("O44").Value = 9
$ws.Range("P44").Value = 1
$ws.Range("Q44").Value = 200
$ws.Range("R44").Value = 20
$ws.Range("S44").Value = 1
$ws.Range("T44").Value = 50
$ws.Range("U44").Value = 0.8360655737704918
$ws.Range("V44").Value = "Data/bombay1.xlsx"
$ws.Range("W44").Value = 517800
$ws.Range("X44").Value = "No es Simulación"

# Row 45
$ws.Range("A45").Value = "2024-04-27 20:05:10"
$ws.Range("B45").Value = 23
$ws.Range("C45").Value = 22
$ws.Range("D45").Value = 3
$ws.Range("E45").Value = 9
$ws.Range("F45").Value = 3
$ws.Range("G45").Value = 7
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0.001
$ws.Range("J45").Value = 0.05
$ws.Range("K45").Value = 0.003
$ws.Range("L45").Value = 100
$ws.Range("M45").Value = 500
$ws.Range("N45").Value = 10
$ws.Range("O45").Value = 9
$ws.Range("P45").Value = 3
$ws.Range("Q45").Value = 200
$ws.Range("R45").Value = 10
$ws.Range("S45").Value = 1
$ws.Range("T45").Value = 90
$ws.Range("U45").Value = 0.9565217391304348
$ws.Range("V45").Value = "Data/bombay1.xlsx"
$ws.Range("W45").Value = 161000
$ws.Range("X45").Value = "No es Simulación"

# Row 46
$ws.Range("A46").Value = "2024-04-27 20:46:17"
$ws.Range("B46").Value = 16
$ws.Range("C46").Value = 14
$ws.Range("D46").Value = 2
$ws.Range("E46").Value = 3
$ws.Range("F46").Value = 4
$ws.Range("G46").Value = 5
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0.001
$ws.Range("J46").Value = 0.05
$ws.Range("K46").Value = 0.003
$ws.Range("L46").Value = 100
$ws.Range("M46").Value = 500
$ws.Range("N46").Value = 10
$ws.Range("O46").Value = 9
$ws.Range("P46").Value = 3
$ws.Range("Q46").Value = 500
$ws.Range("R46").Value = 10
$ws.Range("S46").Value = 1
$ws.Range("T46").Value = 90
$ws.Range("U46").Value = 0.875
$ws.Range("V46").Value = "Data/bombay1.xlsx"
$ws.Range("W46").Value = -376000
$ws.Range("X46").Value = "No es Simulación"

# Row 47
$ws.Range("A47").Value = "2024-04-28 00:07:12"
$ws.Range("B47").Value = 13
$ws.Range("C47").Value = 11
$ws.Range("D47").Value = 2
$ws.Range("E47").Value = 4
$ws.Range("F47").Value = 4
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0.001
$ws.Range("J47").Value = 0.05
$ws.Range("K47").Value = 0.003
$ws.Range("L47").Value = 100
$ws.Range("M47").Value = 500
$ws.Range("N47").Value = 10
$ws.Range("O47").Value = 9
$ws.Range("P47").Value = 3
$ws.Range("Q47").Value = 1000
$ws.Range("R47").Value = 5
$ws.Range("S47").Value = 1
$ws.Range("T47").Value = 90
$ws.Range("U47").Value = 0.8461538461538461
$ws.Range("V47").Value = "Data/bombay1.xlsx"
$ws.Range("W47").Value = 230000
$ws.Range("X47").Value = "No es Simulación"

# Row 48
$ws.Range("A48").Value = "2024-04-28 18:36:21"
$ws.Range("B48").Value = 16
$ws.Range("C48").Value = 11
$ws.Range("D48").Value = 3
$ws.Range("E48").Value = 2
$ws.Range("F48").Value = 3
$ws.Range("G48").Value = 3
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0.001
$ws.Range("J48").Value = 0.05
$ws.Range("K48").Value = 0.003
$ws.Range("L48").Value = 100
$ws.Range("M48").Value = 500
$ws.Range("N48").Value = 10
$ws.Range("O48").Value = 9
$ws.Range("P48").Value = 3
$ws.Range("Q48").Value = 100
$ws.Range("R48").Value = 5
$ws.Range("S48").Value = 1
$ws.Range("T48").Value = 90
$ws.Range("U48").Value = 0.6875
$ws.Range("V48").Value = "Data/bombay1.xlsx"
$ws.Range("W48").Value = -9200
$ws.Range("X48").Value = "No es Simulación"

# Row 49
$ws.Range("A49").Value = "2024-04-29 00:32:25"
$ws.Range("B49").Value = 12
$ws.Range("C49").Value = 10
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 4
$ws.Range("F49").Value = 3
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0.001
$ws.Range("J49").Value = 0.05
$ws.Range("K49").Value = 0.003
$ws.Range("L49").Value = 100
$ws.Range("M49").Value = 500
$ws.Range("N49").Value = 10
$ws.Range("O49").Value = 9
$ws.Range("P49").Value = 3
$ws.Range("Q49").Value = 200
$ws.Range("R49").Value = 5
$ws.Range("S49").Value = 1
$ws.Range("T49").Value = 90
$ws.Range("U49").Value = 0.8333333333333334
$ws.Range("V49").Value = "Data/bombay1.xlsx"
$ws.Range("W49").Value = 15000
$ws.Range("X49").Value = "No es Simulación"
